# PNAD 2009 - "correção nos dados e inicio da analise PNAD 2009"
#
# The original sheet had two section-header rows (row 5 "situação do
# domicílio" and row 8 "grandes regiões e unidades da federação") that
# carried a label but no data, plus a trailing footnote-only row at the
# bottom (row 41). The corrected data set drops those three label-only
# rows entirely, which makes every data row below them shift up by one
# (then two) positions and removes the stray footnote row, shrinking the
# table from A1:H41 down to A1:H38.
#
# Deleting the whole rows (rather than rewriting every value by hand)
# reproduces exactly that shift: Excel moves each row's label together
# with the row, and moves the data values up into the now-vacant rows.
#
# Row 5 ("situação do domicílio") is empty -> delete it; "urbana"/"rural"
# and everything below slide up one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(5).Delete()

# After the first delete, the old row 8 ("grandes regiões e unidades da
# federação", also label-only) is now row 7 -> delete it too; "norte" and
# everything below slide up another row.
$ws.Rows.Item(7).Delete()

# After both deletes, the old trailing footnote row 41 is now row 39 ->
# delete it, dropping the dimension down to A1:H38.
$ws.Rows.Item(39).Delete()

# The column-B sub-header on row 2 is relabelled from the pandas-export
# artifact "unnamed: 1_level_1" to the proper "total".
$ws.Range("B2").Value = "total"
